$d = $word.ActiveDocument

# 1. Update the Sift role description: replace the "(...)" parenthetical
#    with a ": ..." suffix.
$d.Content.Find.Execute(
    "Operations Engineer, Sift, 2009-2011 (Drupal, CentOS, Xen, VMWare/ESXi)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Operations Engineer, Sift, 2009-2011: Drupal, CentOS, Xen, VMWare/ESXi",
    2
)

# 2. Remove the stray "aaaa" paragraph (style FirstParagraph) that
#    followed the Sift role heading.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "aaaa") {
        $p.Range.Delete()
    }
}
